# Update column F (dSF) values on Sheet1 to reflect a repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = -9
    "F3"  = -5
    "F5"  = -5
    "F6"  = 4
    "F9"  = -6
    "F10" = 2
    "F12" = -5
    "F13" = -2
    "F15" = 0
    "F20" = 1
    "F23" = 0
    "F25" = 0
    "F27" = 3
    "F28" = 4
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
